$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $ws.Range("B2").Style
}

$ws.Range("D2").Value = '30.096.22'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '2.118.23'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue $ws "D5" '346.54'
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("E6").Value = '  +0.04%  '
Set-TextValue $ws "D7" '0.5196'
$ws.Range("E7").Value = '  +0.16%  '
Set-TextValue $ws "D8" '0.4477'
$ws.Range("E8").Value = '  +0.24%  '
Set-TextValue $ws "D9" '54.00'
$ws.Range("E9").Value = '  +2.91%  '
Set-TextValue $ws "D10" '0.09372'
$ws.Range("E10").Value = '  -1.50%  '
$ws.Range("E11").Value = '  +0.41%  '
Set-TextValue $ws "D12" '25.42'
$ws.Range("E12").Value = '  +0.05%  '
Set-TextValue $ws "D13" '8.678'
$ws.Range("E13").Value = '  +6.63%  '
Set-TextValue $ws "D14" '6.988'
$ws.Range("E14").Value = '  +3.48%  '
$ws.Range("D15").Value = '2.106.60'
$ws.Range("E15").Value = '  -0.16%  '
$ws.Range("E16").Value = '  +2.94%  '
Set-TextValue $ws "D17" '0.00001170'
$ws.Range("E17").Value = '  -0.16%  '
Set-TextValue $ws "D18" '1.008'
$ws.Range("E18").Value = '  +0.07%  '
Set-TextValue $ws "D19" '21.58'
$ws.Range("E19").Value = '  +4.23%  '
Set-TextValue $ws "D20" '0.06714'
$ws.Range("E20").Value = '  +0.20%  '
Set-TextValue $ws "D21" '6.316'
$ws.Range("E21").Value = '  +1.78%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '30.120.98'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("E24").Value = '  -0.02%  '
Set-TextValue $ws "D25" '2.333'
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("D26").Value = '2.351.65'
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  +0.59%  '
Set-TextValue $ws "D28" '2.554'
$ws.Range("E28").Value = '  +0.27%  '
Set-TextValue $ws "D29" '163.39'
$ws.Range("E29").Value = '  -0.64%  '
Set-TextValue $ws "D30" '134.20'
$ws.Range("E30").Value = '  +0.29%  '
Set-TextValue $ws "D31" '1.157'
$ws.Range("E31").Value = '  -0.61%  '
Set-TextValue $ws "D32" '1.783'
$ws.Range("E32").Value = '  +9.12%  '
Set-TextValue $ws "D33" '0.1059'
$ws.Range("E33").Value = '  +0.14%  '
Set-TextValue $ws "D34" '6.291'
$ws.Range("E34").Value = '  +0.39%  '
Set-TextValue $ws "D35" '6.698'
$ws.Range("E35").Value = '  +8.17%  '
Set-TextValue $ws "D36" '3.974'
$ws.Range("E36").Value = '  +0.75%  '
Set-TextValue $ws "D37" '10.81'
$ws.Range("E37").Value = '  +6.23%  '
Set-TextValue $ws "D38" '0.02647'
$ws.Range("E38").Value = '  +2.55%  '
Set-TextValue $ws "D39" '0.06883'
$ws.Range("E39").Value = '  +1.29%  '
Set-TextValue $ws "D40" '0.7153'
$ws.Range("E40").Value = '  +2.56%  '
Set-TextValue $ws "D41" '12.76'
$ws.Range("E41").Value = '  +1.72%  '
Set-TextValue $ws "D42" '0.2253'
$ws.Range("E42").Value = '  -1.62%  '
$ws.Range("E43").Value = '  +1.32%  '
Set-TextValue $ws "D44" '0.6974'
$ws.Range("E44").Value = '  +3.73%  '
Set-TextValue $ws "D45" '14.69'
$ws.Range("E45").Value = '  +2.94%  '
Set-TextValue $ws "D46" '2.402'
$ws.Range("E47").Value = '  +0.31%  '
Set-TextValue $ws "D48" '3.638'
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("E49").Value = '  +7.41%  '
Set-TextValue $ws "D50" '0.00000000350'
$ws.Range("E50").Value = '  -1.14%  '
$ws.Range("E51").Value = '  +0.42%  '
